$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value2
    if ($v -eq "inactive") {
        $cell.Value = "Inactive"
    }
}

[void]$ws.Cells.Item(70, 2).Select()
